$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.406.33"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.567.01"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "286.73"
$ws.Range("E6").Value = "  +0.36%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3692"
$ws.Range("E7").Value = "  +1.83%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "47.24"
$ws.Range("E8").Value = "  -2.67%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3292"
$ws.Range("E9").Value = "  -1.35%  "
$ws.Range("E10").Value = "  +2.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07437"
$ws.Range("E11").Value = "  +0.56%  "
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.65"
$ws.Range("E13").Value = "  -0.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.899"
$ws.Range("E14").Value = "  -0.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.846"
$ws.Range("E15").Value = "  -0.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.569.64"
$ws.Range("E16").Value = "  +0.26%  "
$ws.Range("E17").Value = "  +0.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06705"
$ws.Range("E18").Value = "  +0.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "86.75"
$ws.Range("E19").Value = "  -1.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.357"
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.37"
$ws.Range("E22").Value = "  +1.34%  "
$ws.Range("E23").Value = "  -1.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.403.10"
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.362"
$ws.Range("E25").Value = "  -1.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.590"
$ws.Range("E26").Value = "  +1.10%  "
$ws.Range("E27").Value = "  +0.26%  "
$ws.Range("E28").Value = "  +0.98%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.935"
$ws.Range("E29").Value = "  -1.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.27"
$ws.Range("E30").Value = "  +0.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.746.01"
$ws.Range("E31").Value = "  +0.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.072"
$ws.Range("E32").Value = "  +1.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.974"
$ws.Range("E33").Value = "  -1.40%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.020"
$ws.Range("E34").Value = "  -1.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.798"
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08261"
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02414"
$ws.Range("E37").Value = "  +0.36%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06341"
$ws.Range("E38").Value = "  -0.59%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.291"
$ws.Range("E39").Value = "  -0.25%  "
$ws.Range("E40").Value = "  -1.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.254"
$ws.Range("E41").Value = "  -1.30%  "
$ws.Range("E42").Value = "  +1.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6128"
$ws.Range("E43").Value = "  +0.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.88"
$ws.Range("E44").Value = "  +1.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5969"
$ws.Range("E45").Value = "  +3.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.753"
$ws.Range("E46").Value = "  -0.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.022"
$ws.Range("E47").Value = "  +0.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.48"
$ws.Range("E48").Value = "  +0.46%  "
$ws.Range("E49").Value = "  -2.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07172"
$ws.Range("E50").Value = "  -0.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "76.54"
$ws.Range("E51").Value = "  +1.27%  "
